$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.021690843904739
$ws.Range("D2").Value = 1.024335558322903
$ws.Range("E2").Value = 1.031020845341375
$ws.Range("F2").Value = 1.038574876731988
$ws.Range("J2").Value = 1.026880887232522
$ws.Range("K2").Value = 1.027164360691542
$ws.Range("L2").Value = 1.033830170034837
$ws.Range("M2").Value = 1.041362520405702
$ws.Range("N2").Value = 1.012902122864228
$ws.Range("C3").Value = 1.022915699364301
$ws.Range("D3").Value = 1.025472823016219
$ws.Range("E3").Value = 1.032148780016428
$ws.Range("F3").Value = 1.039878878316175
$ws.Range("J3").Value = 1.027742038119587
$ws.Range("K3").Value = 1.028107796482225
$ws.Range("L3").Value = 1.034765724350206
$ws.Range("M3").Value = 1.042475259291092
$ws.Range("N3").Value = 1.013198654409765
$ws.Range("C4").Value = 1.02370846351063
$ws.Range("D4").Value = 1.02620923292253
$ws.Range("E4").Value = 1.03287913268215
$ws.Range("F4").Value = 1.040723474215049
$ws.Range("J4").Value = 1.028299001338134
$ws.Range("K4").Value = 1.028718212256756
$ws.Range("L4").Value = 1.035371015438443
$ws.Range("M4").Value = 1.043195552127368
$ws.Range("N4").Value = 1.013390147115894
$ws.Range("C5").Value = 1.024041791797004
$ws.Range("D5").Value = 1.026518946372309
$ws.Range("E5").Value = 1.033186294828395
$ws.Range("F5").Value = 1.041078740419355
$ws.Range("J5").Value = 1.028533087849028
$ws.Range("K5").Value = 1.028974820251849
$ws.Range("L5").Value = 1.035625463062534
$ws.Range("M5").Value = 1.043498430967107
$ws.Range("N5").Value = 1.013470559150397
$ws.Range("C6").Value = 1.024097762092454
$ws.Range("D6").Value = 1.026570956094561
$ws.Range("E6").Value = 1.033237875890924
$ws.Range("F6").Value = 1.041138402848368
$ws.Range("J6").Value = 1.028572388471619
$ws.Range("K6").Value = 1.029017905266288
$ws.Range("L6").Value = 1.035668184991954
$ws.Range("M6").Value = 1.043549289674496
$ws.Range("N6").Value = 1.013484055325069
$ws.Range("C7").Value = 1.023712917259372
$ws.Range("D7").Value = 1.026213370829142
$ws.Range("E7").Value = 1.032883236517172
$ws.Range("F7").Value = 1.040728220518781
$ws.Range("J7").Value = 1.028302129451087
$ws.Range("K7").Value = 1.028721641106223
$ws.Range("L7").Value = 1.035374415444324
$ws.Range("M7").Value = 1.04319959894115
$ws.Range("N7").Value = 1.013391221944717
$ws.Range("C8").Value = 1.022104749281811
$ws.Range("D8").Value = 1.024719794660997
$ws.Range("E8").Value = 1.031401932699479
$ws.Range("F8").Value = 1.039015401568393
$ws.Range("J8").Value = 1.027171971539206
$ws.Range("K8").Value = 1.027483210386509
$ws.Range("L8").Value = 1.034146361289477
$ws.Range("M8").Value = 1.041738519439615
$ws.Range("N8").Value = 1.013002416517993
$ws.Range("C9").Value = 1.019272392115521
$ws.Range("D9").Value = 1.022091871198491
$ws.Range("E9").Value = 1.028795482072915
$ws.Range("F9").Value = 1.036003390439822
$ws.Range("J9").Value = 1.0251784526174
$ws.Range("K9").Value = 1.025300496739728
$ws.Range("L9").Value = 1.031981742428134
$ws.Range("M9").Value = 1.039165957157511
$ws.Range("N9").Value = 1.012314349007886
$ws.Range("C10").Value = 1.017384984856246
$ws.Range("D10").Value = 1.02034248893598
$ws.Range("E10").Value = 1.027060314229686
$ws.Range("F10").Value = 1.03399943006943
$ws.Range("J10").Value = 1.023847995584439
$ws.Range("K10").Value = 1.023844978841961
$ws.Range("L10").Value = 1.030538160295124
$ws.Range("M10").Value = 1.037452191705204
$ws.Range("N10").Value = 1.011853645171917
$ws.Range("C11").Value = 1.016567882665679
$ws.Range("D11").Value = 1.019585575667658
$ws.Range("E11").Value = 1.026309532453585
$ws.Range("F11").Value = 1.033132626194481
$ws.Range("J11").Value = 1.023271533123659
$ws.Range("K11").Value = 1.023214616269054
$ws.Range("L11").Value = 1.029912935901777
$ws.Range("M11").Value = 1.036710392232467
$ws.Range("N11").Value = 1.011653679462174
$ws.Range("C12").Value = 1.016264395459755
$ws.Range("D12").Value = 1.019304509815652
$ws.Range("E12").Value = 1.026030740622867
$ws.Range("F12").Value = 1.032810793105167
$ws.Range("J12").Value = 1.023057352782321
$ws.Range("K12").Value = 1.02298045281515
$ws.Range("L12").Value = 1.029680676484269
$ws.Range("M12").Value = 1.036434893748596
$ws.Range("N12").Value = 1.011579331165185
$ws.Range("C13").Value = 1.016329493592705
$ws.Range("D13").Value = 1.019364795558032
$ws.Range("E13").Value = 1.026090538737243
$ws.Range("F13").Value = 1.032879821258317
$ws.Range("J13").Value = 1.023103297775152
$ws.Range("K13").Value = 1.023030682544738
$ws.Range("L13").Value = 1.02973049799565
$ws.Range("M13").Value = 1.036493987382411
$ws.Range("N13").Value = 1.011595282396544
$ws.Range("C14").Value = 1.016542795907518
$ws.Range("D14").Value = 1.019562340944288
$ws.Range("E14").Value = 1.026286485772477
$ws.Range("F14").Value = 1.033106020600885
$ws.Range("J14").Value = 1.023253830071778
$ws.Range("K14").Value = 1.023195260636443
$ws.Range("L14").Value = 1.029893737749411
$ws.Range("M14").Value = 1.036687618660719
$ws.Range("N14").Value = 1.011647535282064
$ws.Range("C15").Value = 1.016674221194607
$ws.Range("D15").Value = 1.019684066443394
$ws.Range("E15").Value = 1.026407226013614
$ws.Range("F15").Value = 1.033245407492866
$ws.Range("J15").Value = 1.023346570452244
$ws.Range("K15").Value = 1.023296660099797
$ws.Range("L15").Value = 1.029994312020903
$ws.Range("M15").Value = 1.036806926371481
$ws.Range("N15").Value = 1.011679720441768
$ws.Range("C16").Value = 1.017439216323486
$ws.Range("D16").Value = 1.020392734858089
$ws.Range("E16").Value = 1.027110152733824
$ws.Range("F16").Value = 1.03405697625593
$ws.Range("J16").Value = 1.023886245659489
$ws.Range("K16").Value = 1.023886811393775
$ws.Range("L16").Value = 1.030579651209672
$ws.Range("M16").Value = 1.037501428055798
$ws.Range("N16").Value = 1.011866906140499
$ws.Range("C17").Value = 1.017919117811125
$ws.Range("D17").Value = 1.020837418041759
$ws.Range("E17").Value = 1.027551228051581
$ws.Range("F17").Value = 1.034566297647292
$ws.Range("J17").Value = 1.024224670625149
$ws.Range("K17").Value = 1.02425696610334
$ws.Range("L17").Value = 1.030946779363451
$ws.Range("M17").Value = 1.037937142071137
$ws.Range("N17").Value = 1.011984194602441
$ws.Range("C18").Value = 1.018199051623163
$ws.Range("D18").Value = 1.021096850381806
$ws.Range("E18").Value = 1.027808553720032
$ws.Range("F18").Value = 1.034863465530955
$ws.Range("J18").Value = 1.024422032995243
$ws.Range("K18").Value = 1.024472860222667
$ws.Range("L18").Value = 1.031160905240654
$ws.Range("M18").Value = 1.038191313094061
$ws.Range("N18").Value = 1.012052560889342
$ws.Range("C19").Value = 1.018294504512702
$ws.Range("D19").Value = 1.02118531968617
$ws.Range("E19").Value = 1.027896304318693
$ws.Range("F19").Value = 1.03496480732389
$ws.Range("J19").Value = 1.024489322541735
$ws.Range("K19").Value = 1.024546472745361
$ws.Range("L19").Value = 1.031233914345796
$ws.Range("M19").Value = 1.038277983385228
$ws.Range("N19").Value = 1.012075864218767
$ws.Range("C20").Value = 1.017867627352105
$ws.Range("D20").Value = 1.020789701949645
$ws.Range("E20").Value = 1.027503899283301
$ws.Range("F20").Value = 1.034511643060461
$ws.Range("J20").Value = 1.024188364480831
$ws.Range("K20").Value = 1.024217253142233
$ws.Range("L20").Value = 1.030907391423384
$ws.Range("M20").Value = 1.037890391364986
$ws.Range("N20").Value = 1.011971615430321
$ws.Range("C21").Value = 1.016479983136562
$ws.Range("D21").Value = 1.01950416640786
$ws.Range("E21").Value = 1.026228782006593
$ws.Range("F21").Value = 1.03303940678444
$ws.Range("J21").Value = 1.02320950363501
$ws.Range("K21").Value = 1.023146797020599
$ws.Range("L21").Value = 1.029845668365232
$ws.Range("M21").Value = 1.036630598018782
$ws.Range("N21").Value = 1.011632150103402
$ws.Range("C22").Value = 1.015607635870914
$ws.Range("D22").Value = 1.018696392602322
$ws.Range("E22").Value = 1.025427538414887
$ws.Range("F22").Value = 1.032114541764025
$ws.Range("J22").Value = 1.022593727210534
$ws.Range("K22").Value = 1.022473649605804
$ws.Range("L22").Value = 1.02917798567144
$ws.Range("M22").Value = 1.035838739815563
$ws.Range("N22").Value = 1.011418297299844
$ws.Range("C23").Value = 1.016070073231401
$ws.Range("D23").Value = 1.019124562618863
$ws.Range("E23").Value = 1.025852248541194
$ws.Range("F23").Value = 1.032604756279734
$ws.Range("J23").Value = 1.022920193533505
$ws.Range("K23").Value = 1.022830508683171
$ws.Range("L23").Value = 1.029531950268401
$ws.Range("M23").Value = 1.036258498291527
$ws.Range("N23").Value = 1.011531704399516
$ws.Range("C24").Value = 1.017890893626573
$ws.Range("D24").Value = 1.020811262624286
$ws.Range("E24").Value = 1.027525284949973
$ws.Range("F24").Value = 1.034536338839679
$ws.Range("J24").Value = 1.024204769773426
$ws.Range("K24").Value = 1.024235197753533
$ws.Range("L24").Value = 1.030925189180637
$ws.Range("M24").Value = 1.037911515915169
$ws.Range("N24").Value = 1.011977299559188
$ws.Range("C25").Value = 1.020004467820774
$ws.Range("D25").Value = 1.02277079354933
$ws.Range("E25").Value = 1.029468871081094
$ws.Range("F25").Value = 1.036781345557515
$ws.Range("J25").Value = 1.025694074112749
$ws.Range("K25").Value = 1.025864840905051
$ws.Range("L25").Value = 1.032541431143883
$ws.Range("M25").Value = 1.039830793941406
$ws.Range("N25").Value = 1.012492581078627
